# Introduction.docx tidy-up:
#  - fix "user" -> "users" and replace the red TODO note with a forward-looking sentence
#  - rewrite / reorder the "Objectives" bullet list, justify each bullet
#  - drop the old Normal/Configuration/Debug/Autonomous/Supervisor "mode" notes
#    (their useful content now lives in the Objectives bullets)

$d = $word.ActiveDocument

# --- "multiple user to interact" -> "multiple users to interact" (also clears the
#     gramStart/gramEnd proofErr markers that wrapped the old "user") ---
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("multiple user to interact", $true, $false, $false, $false, $false, `
               $true, 1, $false, "multiple users to interact", 2) | Out-Null

# --- Replace the red "Point to the detailed list ... already did it." note with a
#     normal-coloured, forward-looking sentence ---
$p3 = $d.Paragraphs(3)
$p3Text = $p3.Range.Text
$noteIdx = $p3Text.IndexOf(" Point t")
$noteRange = $d.Range($p3.Range.Start + $noteIdx, $p3.Range.End - 1)
$noteRange.Font.Reset()
$noteRange.Text = " A detailed set of objectives for the system is outlined below."

# --- Drop the old "Normal Mode / user mode" ... "Supervisor Mode" paragraphs ---
for ($i = 17; $i -ge 13; $i--) {
    $d.Paragraphs($i).Range.Delete() | Out-Null
}

# --- The old "In autonomous mode..." bullet's content is folded into the
#     "controller" bullet below, so remove this paragraph now ---
$d.Paragraphs(11).Range.Delete() | Out-Null

# --- Bullet: controller control / buggy obstacle behaviour -> becomes the
#     "In autonomous mode..." bullet, with a new "and periodically" insert ---
$p10 = $d.Paragraphs(10)
$p10.Range.ParagraphFormat.Alignment = 3
$r10 = $d.Range($p10.Range.Start, $p10.Range.End - 1)
$r10.Font.Reset()
$r10.Text = "In autonomous mode, the buggy will send"
$r10.InsertAfter(" data back to the connected controller client automatically") | Out-Null
$r10.InsertAfter(" and periodically") | Out-Null
$r10.InsertAfter(". Additionally, the buggy will be able to move around autonomously avoiding obstacles in its way.") | Out-Null

# --- Bullet: environmental data -> becomes the "controller should be able to
#     control ... / buggy should stop moving ..." bullet (red run kept) ---
$p9 = $d.Paragraphs(9)
$p9.Range.ParagraphFormat.Alignment = 3
$r9 = $d.Range($p9.Range.Start, $p9.Range.End - 1)
$r9.Font.Reset()
$r9.Text = "The controller should be able to control the movement of the buggy and configure it if necessary. "
$r9red1 = $d.Range($r9.End, $r9.End)
$r9red1.InsertAfter("The buggy ") | Out-Null
$r9red1.Font.Color = 255
$r9red2 = $d.Range($r9red1.End, $r9red1.End)
$r9red2.InsertAfter("should stop moving if it encounters an object but should let the controller move in another direction to find another path.") | Out-Null
$r9red2.Font.Color = 255

# --- Bullet: reliable/communication -> becomes the environmental-data bullet
#     (unchanged wording) ---
$p8 = $d.Paragraphs(8)
$p8.Range.ParagraphFormat.Alignment = 3
$r8 = $d.Range($p8.Range.Start, $p8.Range.End - 1)
$r8.Font.Reset()
$r8.Text = "The buggy should send back environmental data such as light levels, temperature, and humidity. Any other necessary information should "
$r8.InsertAfter("also be sent back. All the information is to be sent back to the connected controller to be displayed on its GUI.") | Out-Null

# --- Bullet: scalability -> becomes the reliable/communication bullet, extended
#     with new wording about handling lost connections ---
$p7 = $d.Paragraphs(7)
$p7.Range.ParagraphFormat.Alignment = 3
$r7 = $d.Range($p7.Range.Start, $p7.Range.End - 1)
$r7.Font.Reset()
$r7.Text = "The system should "
$r7.InsertAfter("be ") | Out-Null
$r7.InsertAfter("reliable. ") | Out-Null
$r7.InsertAfter("Communication: ") | Out-Null
$r7.InsertAfter("Requests/data should arrive to its destination") | Out-Null
$r7.InsertAfter(" and loss of connection ") | Out-Null
$r7.InsertAfter("between the server, users and buggies should be handled to prevent crashes or unpredictable buggy behaviour.") | Out-Null

# --- Bullet: "system has to be capable ..." stays first, wording reworked ---
$p6 = $d.Paragraphs(6)
$p6.Range.ParagraphFormat.Alignment = 3
$r6 = $d.Range($p6.Range.Start, $p6.Range.End - 1)
$r6.Font.Reset()
$r6.Text = "The system has to be capable of allowing multiple users and multiple buggies "
$r6.InsertAfter("connected to the network ") | Out-Null
$r6.InsertAfter("to interact ") | Out-Null
$r6.InsertAfter("reliably without interference between users or major drop in performance. ") | Out-Null
